$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(2, 1, 2, 140, 132, 3, 1, 0, 5)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(12, $col).Value = $values[$i]
}
